$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff. Numeric-looking text values in column D are
# forced to Text format first so Excel does not silently convert them
# to real numbers (which would drop formatting like trailing zeros).

$ws.Range('D2').Value = '62.298.60'
$ws.Range('E2').Value = '  -3.23%  '
$ws.Range('D3').Value = '3.377.21'
$ws.Range('E3').Value = '  -3.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.48'
$ws.Range('E5').Value = '  -3.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '124.74'
$ws.Range('E6').Value = '  -7.15%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.376.17'
$ws.Range('E8').Value = '  -3.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.472'
$ws.Range('E9').Value = '  -3.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.23'
$ws.Range('E10').Value = '  -5.20%  '
$ws.Range('E11').Value = '  -4.34%  '
$ws.Range('E12').Value = '  -3.91%  '
$ws.Range('D13').Value = '3.952.78'
$ws.Range('E13').Value = '  -3.78%  '
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('D15').Value = '3.381.14'
$ws.Range('E15').Value = '  -3.65%  '
$ws.Range('E16').Value = '  -5.77%  '
$ws.Range('D17').Value = '62.362.13'
$ws.Range('E17').Value = '  -3.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '24.37'
$ws.Range('E18').Value = '  -5.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.21'
$ws.Range('E19').Value = '  -8.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.62'
$ws.Range('E20').Value = '  -2.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.04'
$ws.Range('E21').Value = '  -4.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '370.31'
$ws.Range('E22').Value = '  -6.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.553'
$ws.Range('E23').Value = '  -4.21%  '
$ws.Range('D24').Value = '3.511.95'
$ws.Range('E24').Value = '  -3.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.88'
$ws.Range('E26').Value = '  -5.16%  '
$ws.Range('E27').Value = '  -10.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.85'
$ws.Range('E29').Value = '  -7.26%  '
$ws.Range('E30').Value = '  -6.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.73'
$ws.Range('E31').Value = '  -6.33%  '
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').Value = '3.407.61'
$ws.Range('E33').Value = '  -3.72%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.37'
$ws.Range('E34').Value = '  -6.11%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.147'
$ws.Range('E35').Value = '  -6.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.62'
$ws.Range('E36').Value = '  -3.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.13'
$ws.Range('E37').Value = '  -4.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '164.96'
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('E39').Value = '  -5.60%  '
$ws.Range('E40').Value = '  -5.60%  '
$ws.Range('E41').Value = '  -4.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.764'
$ws.Range('E43').Value = '  -5.75%  '
$ws.Range('E44').Value = '  -2.35%  '
$ws.Range('E45').Value = '  -4.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.53'
$ws.Range('E46').Value = '  -7.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.20'
$ws.Range('E47').Value = '  -11.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.07'
$ws.Range('E48').Value = '  -9.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.57'
$ws.Range('E49').Value = '  -3.33%  '
$ws.Range('D50').Value = '2.228.23'
$ws.Range('E50').Value = '  -6.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.836'
$ws.Range('E51').Value = '  -6.42%  '
